$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column A (rows 1-7) into column B, preserving values
$ws.Range("B1").Value = 40
$ws.Range("B2").Value = 30
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 0.02
$ws.Range("B5").Value = 0.6
$ws.Range("B6").Value = 0.8
$ws.Range("B7").Value = 0.15

# Copy formatting (number format / font) from column A to column B for rows 1-7
$ws.Range("A1:A7").Copy() | Out-Null
$ws.Range("B1:B7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New rows 8 and 9
$ws.Range("A8").Value = -6
$ws.Range("B8").Value = -6
$ws.Range("A9").Value = -2
$ws.Range("B9").Value = -3

# Update selection
$ws.Range("I13").Select() | Out-Null
